$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("FEINmismatch")
$sheet1Dates = @(
    "Mon Feb 24 21:35:59 EST 2025",
    "Mon Feb 24 21:37:17 EST 2025",
    "Mon Feb 24 21:38:33 EST 2025",
    "Mon Feb 24 21:39:47 EST 2025",
    "Mon Feb 24 21:41:01 EST 2025",
    "Mon Feb 24 21:42:15 EST 2025",
    "Mon Feb 24 21:43:29 EST 2025",
    "Mon Feb 24 21:44:43 EST 2025",
    "Mon Feb 24 21:44:54 EST 2025",
    "Mon Feb 24 21:45:07 EST 2025",
    "Mon Feb 24 21:45:19 EST 2025",
    "Mon Feb 24 21:45:31 EST 2025",
    "Mon Feb 24 21:45:43 EST 2025",
    "Mon Feb 24 21:46:56 EST 2025",
    "Mon Feb 24 21:48:10 EST 2025",
    "Mon Feb 24 21:49:23 EST 2025",
    "Mon Feb 24 21:50:37 EST 2025",
    "Mon Feb 24 21:50:49 EST 2025",
    "Mon Feb 24 21:51:01 EST 2025",
    "Mon Feb 24 21:51:12 EST 2025",
    "Mon Feb 24 21:51:24 EST 2025",
    "Mon Feb 24 21:51:36 EST 2025",
    "Mon Feb 24 21:51:47 EST 2025",
    "Mon Feb 24 21:51:59 EST 2025",
    "Mon Feb 24 21:52:10 EST 2025",
    "Mon Feb 24 21:52:22 EST 2025",
    "Mon Feb 24 21:52:34 EST 2025",
    "Mon Feb 24 21:53:48 EST 2025"
)

for ($i = 0; $i -lt $sheet1Dates.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 2).Value = $sheet1Dates[$i]
}

$ws2 = $wb.Worksheets.Item("FEINSSNmismatch")
$sheet2Dates = @(
    "Mon Feb 24 21:54:21 EST 2025",
    "Mon Feb 24 21:54:33 EST 2025",
    "Mon Feb 24 21:54:44 EST 2025",
    "Mon Feb 24 21:54:55 EST 2025",
    "Mon Feb 24 21:55:07 EST 2025",
    "Mon Feb 24 21:55:18 EST 2025",
    "Mon Feb 24 21:55:29 EST 2025",
    "Mon Feb 24 21:55:41 EST 2025",
    "Mon Feb 24 21:55:52 EST 2025",
    "Mon Feb 24 21:56:03 EST 2025",
    "Mon Feb 24 21:56:14 EST 2025",
    "Mon Feb 24 21:56:26 EST 2025",
    "Mon Feb 24 21:56:38 EST 2025",
    "Mon Feb 24 21:56:50 EST 2025",
    "Mon Feb 24 21:57:03 EST 2025",
    "Mon Feb 24 21:57:14 EST 2025",
    "Mon Feb 24 21:57:26 EST 2025",
    "Mon Feb 24 21:57:38 EST 2025",
    "Mon Feb 24 21:57:49 EST 2025",
    "Mon Feb 24 21:58:00 EST 2025",
    "Mon Feb 24 21:58:12 EST 2025",
    "Mon Feb 24 21:58:23 EST 2025",
    "Mon Feb 24 21:58:35 EST 2025",
    "Mon Feb 24 21:58:46 EST 2025"
)

for ($i = 0; $i -lt $sheet2Dates.Length; $i++) {
    $row = $i + 2
    $ws2.Cells.Item($row, 2).Value = $sheet2Dates[$i]
}
